# Apropriações Mauro.xlsx - add a new timesheet entry (row 6) and an empty
# placeholder row (row 7), matching the "Correções e ajustes de acordo com
# as observações levantadas pelo Hermes" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlShiftDown = -4121
$xlFormatFromLeftOrAbove = 0

# --- Row 6: new date/hours entry -------------------------------------------
# Inserting a whole row at position 6 makes Excel inherit the formatting of
# the row above (row 5), so the new cells automatically pick up the same
# date / time-of-day number formats used throughout the sheet.
$ws.Rows(6).Insert($xlShiftDown, $xlFormatFromLeftOrAbove) | Out-Null
$ws.Range("A6").Value = 41527
$ws.Range("B6").Value = 0.027777777777777776

# --- Row 7: empty placeholder row, only column B is present ----------------
$ws.Rows(7).Insert($xlShiftDown, $xlFormatFromLeftOrAbove) | Out-Null
$ws.Range("A7").Clear() | Out-Null

# --- Selection follows the newly added cell, like in the saved workbook ----
$ws.Range("B7").Select() | Out-Null
